$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column P, shifting P:T to Q:U
$ws.Columns("P").Insert()

# Set the new header for column P
$ws.Range("P1").Value = "fin_vyuct_verejne"

# Populate the new fin_vyuct_verejne values for rows 2-151
$ws.Range("P2").Value = 2097983131.007991
$ws.Range("P3").Value = 152584940.45
$ws.Range("P4").Value = 487273834.33
$ws.Range("P5").Value = 2289022251.97
$ws.Range("P6").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("P8").Value = 381472572.3
$ws.Range("P9").Value = 6818677.92
$ws.Range("P10").Value = 0
$ws.Range("P11").Value = 18561765.15
$ws.Range("P12").Value = 775373909.08
$ws.Range("P13").Value = 28725955.31
$ws.Range("P14").Value = 3318112356.43
$ws.Range("P15").Value = 13351311.28
$ws.Range("P16").Value = 1398933631.18
$ws.Range("P17").Value = 958333712.5429
$ws.Range("P18").Value = 356212279.56
$ws.Range("P19").Value = 208299059.8871
$ws.Range("P20").Value = 7874185.189999999
$ws.Range("P21").Value = 4843126.2
$ws.Range("P22").Value = 79415122.48
$ws.Range("P23").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("P26").Value = 758015227.9400001
$ws.Range("P27").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("P29").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("P31").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("P34").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("P36").Value = 5848211.25
$ws.Range("P37").Value = 105612241.96
$ws.Range("P38").Value = 3615900254.8
$ws.Range("P39").Value = 836896913.79
$ws.Range("P40").Value = 113866716.77
$ws.Range("P41").Value = 3643931221.01
$ws.Range("P42").Value = 0
$ws.Range("P43").Value = 0
$ws.Range("P44").Value = 676971589.12
$ws.Range("P45").Value = 5818338927.49
$ws.Range("P46").Value = 176683880.21
$ws.Range("P47").Value = 1163791115.84
$ws.Range("P48").Value = 580227538.5292836
$ws.Range("P49").Value = 355494497.84
$ws.Range("P50").Value = 286535500.4507164
$ws.Range("P51").Value = 3240778538.11
$ws.Range("P52").Value = 447646742.41
$ws.Range("P53").Value = 1568659410.97
$ws.Range("P54").Value = 178006517.92
$ws.Range("P55").Value = 3850585082.325
$ws.Range("P56").Value = 2369283668.3
$ws.Range("P57").Value = 267669865.98
$ws.Range("P58").Value = 9184847364.939142
$ws.Range("P59").Value = 607120319.35
$ws.Range("P60").Value = 289893607.9861853
$ws.Range("P61").Value = 1333346636.059157
$ws.Range("P62").Value = 4141129948.454005
$ws.Range("P63").Value = 1690039619.5117
$ws.Range("P64").Value = 88274441.95999999
$ws.Range("P65").Value = 674642446.25
$ws.Range("P66").Value = 753360246.284
$ws.Range("P67").Value = 2096998391.88
$ws.Range("P68").Value = 24960037.17
$ws.Range("P69").Value = 295610293.2928182
$ws.Range("P70").Value = 580959776.8099999
$ws.Range("P71").Value = 128064328.18
$ws.Range("P72").Value = 4503355865.23
$ws.Range("P73").Value = 149281723.08
$ws.Range("P74").Value = 4844454559.01
$ws.Range("P75").Value = 485077487.76
$ws.Range("P76").Value = 1517814370.17
$ws.Range("P77").Value = 2610339816.13
$ws.Range("P78").Value = 459832271.27
$ws.Range("P79").Value = 0
$ws.Range("P80").Value = 5267096.380000001
$ws.Range("P81").Value = 561662422.8
$ws.Range("P82").Value = 1858459329.56
$ws.Range("P83").Value = 2189559231.1
$ws.Range("P84").Value = 370482997.42
$ws.Range("P85").Value = 13533390.32
$ws.Range("P86").Value = 0
$ws.Range("P87").Value = 0
$ws.Range("P88").Value = 38989183.63
$ws.Range("P89").Value = 0
$ws.Range("P90").Value = 0
$ws.Range("P91").Value = 19193478.71
$ws.Range("P92").Value = 572435491.23
$ws.Range("P93").Value = 0
$ws.Range("P94").Value = 0
$ws.Range("P95").Value = 0
$ws.Range("P96").Value = 0
$ws.Range("P97").Value = 0
$ws.Range("P98").Value = 2213236864.15
$ws.Range("P99").Value = 0
$ws.Range("P100").Value = 0
$ws.Range("P101").Value = 0
$ws.Range("P102").Value = 1542327038.47
$ws.Range("P103").Value = 0
$ws.Range("P104").Value = 0
$ws.Range("P105").Value = 0
$ws.Range("P106").Value = 0
$ws.Range("P107").Value = 0
$ws.Range("P108").Value = 0
$ws.Range("P109").Value = 0
$ws.Range("P110").Value = 0
$ws.Range("P111").Value = 0
$ws.Range("P112").Value = 0
$ws.Range("P113").Value = 0
$ws.Range("P114").Value = 0
$ws.Range("P115").Value = 0
$ws.Range("P116").Value = 0
$ws.Range("P117").Value = 0
$ws.Range("P118").Value = 0
$ws.Range("P119").Value = 0
$ws.Range("P120").Value = 0
$ws.Range("P121").Value = 0
$ws.Range("P122").Value = 0
$ws.Range("P123").Value = 0
$ws.Range("P124").Value = 0
$ws.Range("P125").Value = 0
$ws.Range("P126").Value = 0
$ws.Range("P127").Value = 0
$ws.Range("P128").Value = 0
$ws.Range("P129").Value = 0
$ws.Range("P130").Value = 0
$ws.Range("P131").Value = 226281844.3295787
$ws.Range("P132").Value = 0
$ws.Range("P133").Value = 104692812.01
$ws.Range("P134").Value = 685687781.35
$ws.Range("P135").Value = 1736208313.93
$ws.Range("P136").Value = 0
$ws.Range("P137").Value = 292414422.3192651
$ws.Range("P138").Value = 788781506.34
$ws.Range("P139").Value = 0
$ws.Range("P140").Value = 36627347.54
$ws.Range("P141").Value = 21543025.74
$ws.Range("P142").Value = 90752209.03
$ws.Range("P143").Value = 0
$ws.Range("P144").Value = 42380561.33
$ws.Range("P145").Value = 1399741169.37
$ws.Range("P146").Value = 3546180750.17
$ws.Range("P147").Value = 0
$ws.Range("P148").Value = 0
$ws.Range("P149").Value = 426545296.9611562
$ws.Range("P150").Value = 0
$ws.Range("P151").Value = 0
